$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 108.9481836666667
$ws.Range("H2").Value = 326.844551
$ws.Range("I2").Value = 0.1523660837152667
$ws.Range("J2").Value = 0.1650457680857909
$ws.Range("M2").Value = 9.357204333333334
$ws.Range("N2").Value = 28.071613
$ws.Range("O2").Value = 0.09488819927388772
$ws.Range("P2").Value = 0.09723921246361458
$ws.Range("Q2").Value = 1019.450416314529
$ws.Range("R2").Value = 9175.053746830763
$ws.Range("S2").Value = 0.01445774331415609
$ws.Range("T2").Value = 0.01604892050911469
$ws.Range("G3").Value = 108.9481836666667
$ws.Range("H3").Value = 326.844551
$ws.Range("I3").Value = 0.1523660837152667
$ws.Range("J3").Value = 0.1650457680857909
$ws.Range("O3").Value = 0.03515788078400975
$ws.Range("P3").Value = 0.03602897584196811
$ws.Range("Q3").Value = 377.7257496323741
$ws.Range("R3").Value = 3399.531746691366
$ws.Range("S3").Value = 0.005356868606787796
$ws.Range("T3").Value = 0.005946429991182034
$ws.Range("G4").Value = 108.9481836666667
$ws.Range("H4").Value = 326.844551
$ws.Range("I4").Value = 0.1523660837152667
$ws.Range("J4").Value = 0.1650457680857909
$ws.Range("M4").Value = 42.666574
$ws.Range("N4").Value = 127.999722
$ws.Range("O4").Value = 0.4326670907061247
$ws.Range("P4").Value = 0.4433871385602816
$ws.Range("Q4").Value = 4648.445740579426
$ws.Range("R4").Value = 41836.01166521483
$ws.Range("S4").Value = 0.0659237901633703
$ws.Range("T4").Value = 0.07317917084304269
$ws.Range("G5").Value = 108.9481836666667
$ws.Range("H5").Value = 326.844551
$ws.Range("I5").Value = 0.1523660837152667
$ws.Range("J5").Value = 0.1650457680857909
$ws.Range("M5").Value = 7.15268
$ws.Range("N5").Value = 14.30536
$ws.Range("O5").Value = 0.07253287424370854
$ws.Range("P5").Value = 0.04955333134610019
$ws.Range("Q5").Value = 779.2714943488934
$ws.Range("R5").Value = 4675.628966093361
$ws.Range("S5").Value = 0.01105154998912581
$ws.Range("T5").Value = 0.008178567633226807
$ws.Range("G6").Value = 108.9481836666667
$ws.Range("H6").Value = 326.844551
$ws.Range("I6").Value = 0.1523660837152667
$ws.Range("J6").Value = 0.1650457680857909
$ws.Range("M6").Value = 35.96946000000001
$ws.Range("N6").Value = 107.90838
$ws.Range("O6").Value = 0.3647539549922693
$ws.Range("P6").Value = 0.3737913417880355
$ws.Range("Q6").Value = 3918.807334470821
$ws.Range("R6").Value = 35269.26601023738
$ws.Range("S6").Value = 0.05557613164182672
$ws.Range("T6").Value = 0.06169267910922472
$ws.Range("I7").Value = 0.2954065074566193
$ws.Range("J7").Value = 0.3199898083081954
$ws.Range("M7").Value = 9.357204333333334
$ws.Range("N7").Value = 28.071613
$ws.Range("O7").Value = 0.09488819927388772
$ws.Range("P7").Value = 0.09723921246361458
$ws.Range("Q7").Value = 1976.504742167215
$ws.Range("R7").Value = 17788.54267950494
$ws.Range("S7").Value = 0.02803059154634689
$ws.Range("T7").Value = 0.03111555695627191
$ws.Range("I8").Value = 0.2954065074566193
$ws.Range("J8").Value = 0.3199898083081954
$ws.Range("O8").Value = 0.03515788078400975
$ws.Range("P8").Value = 0.03602897584196811
$ws.Range("S8").Value = 0.01038586677198051
$ws.Range("T8").Value = 0.01152890507321198
$ws.Range("I9").Value = 0.2954065074566193
$ws.Range("J9").Value = 0.3199898083081954
$ws.Range("M9").Value = 42.666574
$ws.Range("N9").Value = 127.999722
$ws.Range("O9").Value = 0.4326670907061247
$ws.Range("P9").Value = 0.4433871385602816
$ws.Range("Q9").Value = 9012.380497304706
$ws.Range("R9").Value = 81111.42447574236
$ws.Range("S9").Value = 0.1278126741569126
$ws.Range("T9").Value = 0.1418793654742238
$ws.Range("I10").Value = 0.2954065074566193
$ws.Range("J10").Value = 0.3199898083081954
$ws.Range("M10").Value = 7.15268
$ws.Range("N10").Value = 14.30536
$ws.Range("O10").Value = 0.07253287424370854
$ws.Range("P10").Value = 0.04955333134610019
$ws.Range("Q10").Value = 1510.847197046133
$ws.Range("R10").Value = 9065.0831822768
$ws.Range("S10").Value = 0.02142668305612411
$ws.Range("T10").Value = 0.01585656099847109
$ws.Range("I11").Value = 0.2954065074566193
$ws.Range("J11").Value = 0.3199898083081954
$ws.Range("M11").Value = 35.96946000000001
$ws.Range("N11").Value = 107.90838
$ws.Range("O11").Value = 0.3647539549922693
$ws.Range("P11").Value = 0.3737913417880355
$ws.Range("Q11").Value = 7597.7616530116
$ws.Range("R11").Value = 68379.8548771044
$ws.Range("S11").Value = 0.1077506919252552
$ws.Range("T11").Value = 0.1196094198060166
$ws.Range("G12").Value = 109.1710686666667
$ws.Range("H12").Value = 327.513206
$ws.Range("I12").Value = 0.1526777925792968
$ws.Range("J12").Value = 0.1653834169091284
$ws.Range("M12").Value = 9.357204333333334
$ws.Range("N12").Value = 28.071613
$ws.Range("O12").Value = 0.09488819927388772
$ws.Range("P12").Value = 0.09723921246361458
$ws.Range("Q12").Value = 1021.535996802364
$ws.Range("R12").Value = 9193.823971221278
$ws.Range("S12").Value = 0.01448732080696161
$ws.Range("T12").Value = 0.01608175321478528
$ws.Range("G13").Value = 109.1710686666667
$ws.Range("H13").Value = 327.513206
$ws.Range("I13").Value = 0.1526777925792968
$ws.Range("J13").Value = 0.1653834169091284
$ws.Range("O13").Value = 0.03515788078400975
$ws.Range("P13").Value = 0.03602897584196811
$ws.Range("Q13").Value = 378.498496830844
$ws.Range("R13").Value = 3406.486471477596
$ws.Range("S13").Value = 0.005367827629868685
$ws.Range("T13").Value = 0.005958595132481127
$ws.Range("G14").Value = 109.1710686666667
$ws.Range("H14").Value = 327.513206
$ws.Range("I14").Value = 0.1526777925792968
$ws.Range("J14").Value = 0.1653834169091284
$ws.Range("M14").Value = 42.666574
$ws.Range("N14").Value = 127.999722
$ws.Range("O14").Value = 0.4326670907061247
$ws.Range("P14").Value = 0.4433871385602816
$ws.Range("Q14").Value = 4657.955479925416
$ws.Range("R14").Value = 41921.59931932874
$ws.Range("S14").Value = 0.06605865633071752
$ws.Range("T14").Value = 0.07332887998866053
$ws.Range("G15").Value = 109.1710686666667
$ws.Range("H15").Value = 327.513206
$ws.Range("I15").Value = 0.1526777925792968
$ws.Range("J15").Value = 0.1653834169091284
$ws.Range("M15").Value = 7.15268
$ws.Range("N15").Value = 14.30536
$ws.Range("O15").Value = 0.07253287424370854
$ws.Range("P15").Value = 0.04955333134610019
$ws.Range("Q15").Value = 780.8657194306934
$ws.Range("R15").Value = 4685.19431658416
$ws.Range("S15").Value = 0.01107415912896115
$ws.Range("T15").Value = 0.008195299257248267
$ws.Range("G16").Value = 109.1710686666667
$ws.Range("H16").Value = 327.513206
$ws.Range("I16").Value = 0.1526777925792968
$ws.Range("J16").Value = 0.1653834169091284
$ws.Range("M16").Value = 35.96946000000001
$ws.Range("N16").Value = 107.90838
$ws.Range("O16").Value = 0.3647539549922693
$ws.Range("P16").Value = 0.3737913417880355
$ws.Range("Q16").Value = 3926.82438756292
$ws.Range("R16").Value = 35341.41948806628
$ws.Range("S16").Value = 0.05568982868278785
$ws.Range("T16").Value = 0.06181888931595318
$ws.Range("G17").Value = 164.799919
$ws.Range("H17").Value = 329.599838
$ws.Range("I17").Value = 0.2304757859153342
$ws.Range("J17").Value = 0.166437097565877
$ws.Range("M17").Value = 9.357204333333334
$ws.Range("N17").Value = 28.071613
$ws.Range("O17").Value = 0.09488819927388772
$ws.Range("P17").Value = 0.09723921246361458
$ws.Range("Q17").Value = 1542.066516199782
$ws.Range("R17").Value = 9252.399097198693
$ws.Range("S17").Value = 0.02186943230174011
$ws.Range("T17").Value = 0.01618421229203566
$ws.Range("G18").Value = 164.799919
$ws.Range("H18").Value = 329.599838
$ws.Range("I18").Value = 0.2304757859153342
$ws.Range("J18").Value = 0.166437097565877
$ws.Range("O18").Value = 0.03515788078400975
$ws.Range("P18").Value = 0.03602897584196811
$ws.Range("Q18").Value = 571.364944771218
$ws.Range("R18").Value = 3428.189668627308
$ws.Range("S18").Value = 0.008103040204812272
$ws.Range("T18").Value = 0.005996558167408271
$ws.Range("G19").Value = 164.799919
$ws.Range("H19").Value = 329.599838
$ws.Range("I19").Value = 0.2304757859153342
$ws.Range("J19").Value = 0.166437097565877
$ws.Range("M19").Value = 42.666574
$ws.Range("N19").Value = 127.999722
$ws.Range("O19").Value = 0.4326670907061247
$ws.Range("P19").Value = 0.4433871385602816
$ws.Range("Q19").Value = 7031.447939207506
$ws.Range("R19").Value = 42188.68763524503
$ws.Range("S19").Value = 0.09971928777019527
$ws.Range("T19").Value = 0.07379606844001262
$ws.Range("G20").Value = 164.799919
$ws.Range("H20").Value = 329.599838
$ws.Range("I20").Value = 0.2304757859153342
$ws.Range("J20").Value = 0.166437097565877
$ws.Range("M20").Value = 7.15268
$ws.Range("N20").Value = 14.30536
$ws.Range("O20").Value = 0.07253287424370854
$ws.Range("P20").Value = 0.04955333134610019
$ws.Range("Q20").Value = 1178.76108463292
$ws.Range("R20").Value = 4715.044338531679
$ws.Range("S20").Value = 0.01671707119601682
$ws.Range("T20").Value = 0.008247512643965107
$ws.Range("G21").Value = 164.799919
$ws.Range("H21").Value = 329.599838
$ws.Range("I21").Value = 0.2304757859153342
$ws.Range("J21").Value = 0.166437097565877
$ws.Range("M21").Value = 35.96946000000001
$ws.Range("N21").Value = 107.90838
$ws.Range("O21").Value = 0.3647539549922693
$ws.Range("P21").Value = 0.3737913417880355
$ws.Range("Q21").Value = 5927.76409447374
$ws.Range("R21").Value = 35566.58456684244
$ws.Range("S21").Value = 0.08406695444256967
$ws.Range("T21").Value = 0.06221274602245534
$ws.Range("G22").Value = 120.894928
$ws.Range("H22").Value = 362.684784
$ws.Range("I22").Value = 0.169073830333483
$ws.Range("J22").Value = 0.1831439091310082
$ws.Range("M22").Value = 9.357204333333334
$ws.Range("N22").Value = 28.071613
$ws.Range("O22").Value = 0.09488819927388772
$ws.Range("P22").Value = 0.09723921246361458
$ws.Range("Q22").Value = 1131.238544159621
$ws.Range("R22").Value = 10181.14689743659
$ws.Range("S22").Value = 0.01604311130468302
$ws.Range("T22").Value = 0.01780876949140703
$ws.Range("G23").Value = 120.894928
$ws.Range("H23").Value = 362.684784
$ws.Range("I23").Value = 0.169073830333483
$ws.Range("J23").Value = 0.1831439091310082
$ws.Range("O23").Value = 0.03515788078400975
$ws.Range("P23").Value = 0.03602897584196811
$ws.Range("Q23").Value = 419.145375064416
$ws.Range("R23").Value = 3772.308375579744
$ws.Range("S23").Value = 0.005944277570560486
$ws.Range("T23").Value = 0.006598487477684698
$ws.Range("G24").Value = 120.894928
$ws.Range("H24").Value = 362.684784
$ws.Range("I24").Value = 0.169073830333483
$ws.Range("J24").Value = 0.1831439091310082
$ws.Range("M24").Value = 42.666574
$ws.Range("N24").Value = 127.999722
$ws.Range("O24").Value = 0.4326670907061247
$ws.Range("P24").Value = 0.4433871385602816
$ws.Range("Q24").Value = 5158.172391736673
$ws.Range("R24").Value = 46423.55152563006
$ws.Range("S24").Value = 0.07315268228492904
$ws.Range("T24").Value = 0.08120365381434197
$ws.Range("G25").Value = 120.894928
$ws.Range("H25").Value = 362.684784
$ws.Range("I25").Value = 0.169073830333483
$ws.Range("J25").Value = 0.1831439091310082
$ws.Range("M25").Value = 7.15268
$ws.Range("N25").Value = 14.30536
$ws.Range("O25").Value = 0.07253287424370854
$ws.Range("P25").Value = 0.04955333134610019
$ws.Range("Q25").Value = 864.72273360704
$ws.Range("R25").Value = 5188.336401642241
$ws.Range("S25").Value = 0.01226341087348064
$ws.Range("T25").Value = 0.009075390813188915
$ws.Range("G26").Value = 120.894928
$ws.Range("H26").Value = 362.684784
$ws.Range("I26").Value = 0.169073830333483
$ws.Range("J26").Value = 0.1831439091310082
$ws.Range("M26").Value = 35.96946000000001
$ws.Range("N26").Value = 107.90838
$ws.Range("O26").Value = 0.3647539549922693
$ws.Range("P26").Value = 0.3737913417880355
$ws.Range("Q26").Value = 4348.525276898881
$ws.Range("R26").Value = 39136.72749208993
$ws.Range("S26").Value = 0.06167034829982983
$ws.Range("T26").Value = 0.06845760753438561
